$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FX OHLCV rows to append (dates as Excel serial numbers, matching
# the existing column A datetime serials already present in the sheet).
$newRows = @(
    @(791, 45142, 0.007012,  0.007063,  0.007001,  0.007053, 0),
    @(792, 45145, 0.007053,  0.007062,  0.007014,  0.007016, 0),
    @(793, 45146, 0.007016,  0.007021,  0.00697,   0.006971, 0),
    @(794, 45147, 0.006974,  0.006992,  0.006954,  0.006957, 0),
    @(795, 45148, 0.006957,  0.006974,  0.006906,  0.006908, 0),
    @(796, 45149, 0.006907,  0.006923,  0.006895,  0.006895, 0),
    @(797, 45152, 0.006897,  0.006908,  0.006867,  0.006867, 0),
    @(798, 45153, 0.006869,  0.006891,  0.006856,  0.006866, 0),
    @(799, 45154, 0.006867,  0.006881,  0.006831,  0.006832, 0),
    @(800, 45155, 0.006833,  0.006866,  0.006826,  0.006853, 0),
    @(801, 45156, 0.006856,  0.006898,  0.006856,  0.006874, 0),
    @(802, 45159, 0.006874,  0.0068884, 0.006832,  0.006837, 0),
    @(803, 45160, 0.006837,  0.006872,  0.006832,  0.006851, 0),
    @(804, 45161, 0.0068515, 0.0068682, 0.0068515, 0.006861, 0)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
    $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
    $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
    $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
    $ws.Cells.Item($rowNum, 5).Value2 = $r[5]
    $ws.Cells.Item($rowNum, 6).Value2 = $r[6]
}

# Match the styling of the datetime column (A) used by the existing data
# rows (centered/top aligned, bordered, custom datetime number format).
$ws.Cells.Item(790, 1).Copy()
$ws.Range("A791:A804").PasteSpecial(-4122)
$excel.CutCopyMode = $false
